$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.230.29'
$ws.Range("E2").Value = '  +0.27%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.842.12'
$ws.Range("E3").Value = '  +0.29%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9991'
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.58'
$ws.Range("E5").Value = '  +0.08%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6742'
$ws.Range("E6").Value = '  -1.67%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07420'
$ws.Range("E8").Value = '  -0.54%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2951'
$ws.Range("E9").Value = '  -1.99%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.82'
$ws.Range("E10").Value = '  -1.38%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07718'
$ws.Range("E11").Value = '  +0.70%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.857.57'
$ws.Range("E12").Value = '  +1.09%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6711'
$ws.Range("E14").Value = '  -1.57%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '86.01'
$ws.Range("E15").Value = '  -1.78%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.117'
$ws.Range("E16").Value = '  -0.59%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.300.74'
$ws.Range("E17").Value = '  +0.55%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008310'
$ws.Range("E18").Value = '  +1.68%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '228.44'
$ws.Range("E19").Value = '  +0.28%  '

$ws.Range("E20").Value = '  -0.28%  '

$ws.Range("E21").Value = '  +0.17%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.198'
$ws.Range("E22").Value = '  -2.64%  '

$ws.Range("E23").Value = '  +0.04%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '160.74'
$ws.Range("E24").Value = '  +0.40%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.693'
$ws.Range("E25").Value = '  -0.74%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1401'
$ws.Range("E26").Value = '  -3.66%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.511'
$ws.Range("E28").Value = '  -0.04%  '

$ws.Range("E29").Value = '  -2.32%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.067'
$ws.Range("E30").Value = '  -2.00%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.191'
$ws.Range("E31").Value = '  -0.71%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05312'
$ws.Range("E32").Value = '  +2.61%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7621'
$ws.Range("E33").Value = '  -0.51%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.869'
$ws.Range("E34").Value = '  +1.40%  '

$ws.Range("E35").Value = '  -0.02%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.677'
$ws.Range("E36").Value = '  +0.11%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.328.18'
$ws.Range("E37").Value = '  +1.19%  '

$ws.Range("E38").Value = '  -1.59%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.721'
$ws.Range("E39").Value = '  +0.05%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9205'
$ws.Range("E40").Value = '  -1.51%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.940'
$ws.Range("E41").Value = '  +2.52%  '

$ws.Range("E42").Value = '  +0.21%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '103.55'
$ws.Range("E43").Value = '  -1.08%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.08197'
$ws.Range("E44").Value = '  +15.82%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5161'
$ws.Range("E47").Value = '  -0.79%  '

$ws.Range("E48").Value = '  +0.25%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '63.79'
$ws.Range("E49").Value = '  -1.89%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.143'
$ws.Range("E50").Value = '  -3.88%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05957'
$ws.Range("E51").Value = '  +0.34%  '

# Row 45/46: RocketPoolETH and BabyDogeCoin swap positions with updated data
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000126'
$ws.Range("E45").Value = '  +2.29%  '

$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.002.28'
$ws.Range("E46").Value = '  +1.00%  '
